$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - update existing component (Phoenix Contact header -> CUI Devices 30A header)
$ws.Range("A6").Value = "1x2 headers 30 A"
$ws.Range("B6").Value = 1.54
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = "https://www.mouser.be/ProductDetail/Phoenix-Contact/1714955?qs=sGAEpiMZZMvZTcaMAxB2AKJ8wpDgQJg56QuuQcQ95jA%3D"

# Row 7 - new component 2x15
$ws.Range("A7").Value = "2x15"
$ws.Range("C7").Value = 1
$ws.Range("E7").Value = "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-2-aantal-polen-per-rij-15-blg2x15-1-stuks-1492273"

# Row 8 - new component 2x17
$ws.Range("A8").Value = "2x17"
$ws.Range("E8").Value = "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-2-aantal-polen-per-rij-17-blg2x17-1-stuks-1492276"

# Row 9 - new component 2x10
$ws.Range("A9").Value = "2x10"
$ws.Range("E9").Value = "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-2-aantal-polen-per-rij-10-blg2x10-1-stuks-1492267"

# Row 10 - new component 2x8
$ws.Range("A10").Value = "2x8"
$ws.Range("E10").Value = "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-2-aantal-polen-per-rij-8-blg2x8-1-stuks-1492298"

# Row 11 - new component 1x8
$ws.Range("A11").Value = "1x8"
$ws.Range("E11").Value = "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-1-aantal-polen-per-rij-8-blg1x8-1-stuks-1492302"

# Row 12 - new component 1x2 headers 15 A
$ws.Range("A12").Value = "1x2 headers 15 A"
$ws.Range("E12").Value = "https://www.mouser.be/ProductDetail/CUI-Devices/TB002-500-02BE?qs=sGAEpiMZZMvZTcaMAxB2AHpdXjUJWjdtGYWJDK8ID%2FsZJpc5bbOw%2FQ%3D%3D"

# Now convert the relevant URL cells into actual hyperlinks (matching diff's <hyperlinks> list)
# Order matches the target rId1..rId4 allocation order.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E10"), "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-2-aantal-polen-per-rij-8-blg2x8-1-stuks-1492298") | Out-Null
$ws.Range("E10").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.conrad.be/p/econ-connect-female-header-standaard-aantal-rijen-1-aantal-polen-per-rij-8-blg1x8-1-stuks-1492302") | Out-Null
$ws.Range("E11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.mouser.be/ProductDetail/Phoenix-Contact/1714955?qs=sGAEpiMZZMvZTcaMAxB2AKJ8wpDgQJg56QuuQcQ95jA%3D") | Out-Null
$ws.Range("E6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E12"), "https://www.mouser.be/ProductDetail/CUI-Devices/TB002-500-02BE?qs=sGAEpiMZZMvZTcaMAxB2AHpdXjUJWjdtGYWJDK8ID%2FsZJpc5bbOw%2FQ%3D%3D") | Out-Null
$ws.Range("E12").Style = "Hyperlink"

$wb.Save()
